# Regenerate save_data: recompute column G ("K") values for each row
# (commit: "regen save_data to use K instead of Strike#, regen std/mean,
#  calc and write s_vals")

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New K values, row-indexed (row 2 .. row 24), matching the recalculated
# s_vals written by the regenerated pipeline.
$newK = @{
    2  = 10
    3  = 8
    4  = 7
    5  = 10
    6  = 4
    7  = 2
    8  = 8
    9  = 6
    10 = 2
    11 = 11
    12 = 4
    13 = 4
    14 = 9
    15 = 10
    16 = 6
    17 = 6
    18 = 5
    19 = 4
    20 = 6
    21 = 7
    22 = 4
    23 = 7
    24 = 4
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
